$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the Abstract (D13) and Authors (E13) cells on row 13 to their
# earlier values, undoing the bug-fix/dump-file commit this reverts.
$abstract13 = "Background`nid=`"Par1`">This study aimed to assess the impact of 15" + ([char]0xA0) + "days before, 15" + ([char]0xA0) + "days during, and 15" + ([char]0xA0) + "days after the lockdown on the trends in the prevalence and mortality in 27 countries during COVID-19 pandemic.`n`n`nMethods`nid=`"Par2`">Twenty-seven countries were randomly selected from the different continents.`n`n The information on the trends in the prevalence and mortality due to COVID-19 pandemic in 27 countries was obtained from World Health Organization and lockdown data were obtained from concerned countries and their ministries.`n`n The impact of lockdown for 15" + ([char]0xA0) + "days before, 15" + ([char]0xA0) + "days during, and 15" + ([char]0xA0) + "days after the lockdown on the prevalence and mortality due to the COVID-19 pandemic in 27 countries was analyzed.`n`n`nResults`nid=`"Par3`">The findings showed that 15" + ([char]0xA0) + "days after the lockdown there was a trend toward a decline, but no significant decline in the mean prevalence and mean mortality rate due to the COVID-19 pandemic compared to 15" + ([char]0xA0) + "days before, and 15" + ([char]0xA0) + "days during the lockdown in 27 countries.`n`n The mean growth factor for number of cases was 1.18 and for mortality rate was 1.16.`nConclusions`nid=`"Par4`">The findings indicate that 15" + ([char]0xA0) + "days after the lockdown, daily cases of COVID-19 and the growth factor of the disease showed a declined trend, but there was no significant decline in the prevalence and mortality.`n`n`n"
$authors13 = "[Sultan Ayoub%Meo%sultanmeo@hotmail.com%1,   Abdulelah Adnan%Abukhalaf%Abdulelahabukhalaf@gmail.com%1,   Ali Abdullah%Alomar%AliAlomarMD@gmail.com%1,   Faris Jamal%AlMutairi%faris11300@gmail.com%1,   Adnan Mahmood%Usmani%adnan_mahmood71@hotmail.com%1,   David C.%Klonoff%dklonoff@diabetestechnology.org%1]"

$ws.Range("D13").Value = $abstract13
$ws.Range("E13").Value = $authors13
